$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 496
$ws1.Range("F6").Value = 933
$ws1.Range("F7").Value = 173
$ws1.Range("F8").Value = 3
$ws1.Range("F9").Value = 970
$ws1.Range("F10").Value = 759
$ws1.Range("F11").Value = 208
$ws1.Range("F13").Value = 80
$ws1.Range("F14").Value = 793
$ws1.Range("F15").Value = 260
$ws1.Range("F16").Value = 561
$ws1.Range("F17").Value = 494
$ws1.Range("F21").Value = 1128
$ws1.Range("F22").Value = 2817
$ws1.Range("F23").Value = 1335
$ws1.Range("F24").Value = 668
$ws1.Range("F25").Value = 172
$ws1.Range("F28").Value = 982
$ws1.Range("F30").Value = 1872
$ws1.Range("F31").Value = 35
$ws1.Range("F32").Value = 5
$ws1.Range("F33").Value = 1349

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 723

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 723
$ws4.Range("F7").Value = 496
$ws4.Range("F13").Value = 933
$ws4.Range("F14").Value = 173
$ws4.Range("F16").Value = 3
$ws4.Range("F17").Value = 970
$ws4.Range("F18").Value = 759
$ws4.Range("F19").Value = 208
$ws4.Range("F25").Value = 80
$ws4.Range("F27").Value = 793
$ws4.Range("F28").Value = 260
$ws4.Range("F29").Value = 561
$ws4.Range("F30").Value = 494
$ws4.Range("F34").Value = 1128
$ws4.Range("F35").Value = 2817
$ws4.Range("F36").Value = 1335
$ws4.Range("F37").Value = 668
$ws4.Range("F38").Value = 172
$ws4.Range("F43").Value = 982
$ws4.Range("F45").Value = 1872
$ws4.Range("F46").Value = 35
$ws4.Range("F47").Value = 5
$ws4.Range("F48").Value = 1349
